# Add "Wins" / "Losses" / "Ties" columns (AD, AE, AF) to the season-record
# worksheet. Column AC ("Unnamed: 28") is the last existing column, so the
# new columns are appended directly after it.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---------------------------------------------------
# Copy the formatting of an existing header cell (bold font, thin border,
# centered/top alignment - style index 1) onto each new header cell before
# setting its text so the new headers look like the rest of the header row.
$ws.Range("AC1").Copy($ws.Range("AD1"))
$ws.Range("AD1").Value = "Wins"

$ws.Range("AC1").Copy($ws.Range("AE1"))
$ws.Range("AE1").Value = "Losses"

$ws.Range("AC1").Copy($ws.Range("AF1"))
$ws.Range("AF1").Value = "Ties"

# --- Data rows (rows 2-43) -------------------------------------------------
# Every team/player row gets the same season record: 80 wins, 82 losses, 0 ties.
for ($r = 2; $r -le 43; $r++) {
    $ws.Cells.Item($r, 30).Value = 80
    $ws.Cells.Item($r, 31).Value = 82
    $ws.Cells.Item($r, 32).Value = 0
}

Write-Output "Added Wins/Losses/Ties columns (AD:AF) for rows 1-43"
